# Generate Report for Handoff
# Refresh the "Latest Handoff" timestamps for the 8666816a-aab3-472d-bea9-9e0ef1f71c72
# file (row 6 on each sheet) to reflect a newly regenerated handoff xliff.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 8666816a-aab3-472d-bea9-9e0ef1f71c72.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-19 02:39:30"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the same file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-19 02:39:25"

# de-de sheet: "Latest Handoff Datetime" column (H) for the same file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-19 02:39:30"
